# "Revert to 2.1.1 files" - update the passenger (SoCDTtiNTY-psgr) sheet's
# Share-of-New-This-Year assumptions and leave the active tab on that sheet.

$wb = $excel.ActiveWorkbook
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# LDVs (row 2): B2 becomes a calibration formula, D2 gets a lower share.
$wsPsgr.Range("B2").Formula = "=0.076+(0.076-0.0725)"
$wsPsgr.Range("D2").Value = 0.0735

# motorbikes (row 5): B5 and E5 drop to 0.01.
$wsPsgr.Range("B5").Value = 0.01
$wsPsgr.Range("E5").Value = 0.01

# Make the psgr sheet the active/selected tab, with E6 as the selection,
# matching the saved workbook view state.
$wsPsgr.Activate()
$wsPsgr.Range("E6").Select()
